$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E ("modelscope") to hold the new
# "模型资源" (model resource) field. Excel will shift all subsequent
# columns (old E..U -> F..V) and auto-update formulas/col widths.
$ws.Columns("E").Insert()

# Header block describing the new field (rows 1-6 are header/meta rows).
$ws.Range("E2").Value = "模型资源文件名称"
$ws.Range("E3").Value = "modelresource"
$ws.Range("E4").Value = "模型资源"
$ws.Range("E6").Value = "string"

# Data rows 7-13: populate the new column with the resource file name.
$ws.Range("E7").Value = "puluomixiusi"
$ws.Range("E8").Value = "puluomixiusi"
$ws.Range("E9").Value = "puluomixiusi"
$ws.Range("E10").Value = "puluomixiusi"
$ws.Range("E11").Value = "puluomixiusi"
$ws.Range("E12").Value = "puluomixiusi"
$ws.Range("E13").Value = "puluomixiusi"

# Match the final saved selection state.
$ws.Range("E9").Select()
